$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '56.752.96'
$ws.Range("E2").Value = '  +10.82%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.256.11'
$ws.Range("E3").Value = '  +6.23%  '

# Row 4
$ws.Range("E4").Value = '  +0.06%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '398.21'
$ws.Range("E5").Value = '  +2.71%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '111.39'
$ws.Range("E6").Value = '  +9.04%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.561'
$ws.Range("E7").Value = '  +4.66%  '

# Row 8
$ws.Range("E8").Value = '  -0.09%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.621'
$ws.Range("E9").Value = '  +7.31%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.49'
$ws.Range("E10").Value = '  +7.37%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0956'
$ws.Range("E11").Value = '  +12.64%  '

# Row 12
$ws.Range("E12").Value = '  +2.33%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.743.75'
$ws.Range("E13").Value = '  +5.62%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '19.31'
$ws.Range("E14").Value = '  +5.52%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.12'
$ws.Range("E15").Value = '  +5.59%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.261.91'
$ws.Range("E16").Value = '  +6.73%  '

# Row 17
$ws.Range("E17").Value = '  +5.81%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.14'
$ws.Range("E18").Value = '  +4.77%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '56.558.83'
$ws.Range("E19").Value = '  +10.50%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.33'
$ws.Range("E20").Value = '  +3.78%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0000105'
$ws.Range("E21").Value = '  +9.52%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '13.03'
$ws.Range("E22").Value = '  +6.06%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '300.16'
$ws.Range("E23").Value = '  +13.49%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '75.59'
$ws.Range("E24").Value = '  +8.50%  '

# Row 25
$ws.Range("E25").Value = '  +3.64%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.13'
$ws.Range("E26").Value = '  +2.73%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '28.39'
$ws.Range("E27").Value = '  +4.79%  '

# Row 28
$ws.Range("E28").Value = '  +4.33%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.37'
$ws.Range("E29").Value = '  +1.77%  '

# Row 30
$ws.Range("E30").Value = '  +4.50%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("E31").Value = '  -0.17%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.112'
$ws.Range("E32").Value = '  +6.67%  '

# Row 33
$ws.Range("E33").Value = '  +6.80%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '36.75'
$ws.Range("E34").Value = '  +3.00%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0486'
$ws.Range("E35").Value = '  +3.27%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.12'
$ws.Range("E36").Value = '  +2.20%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '51.71'
$ws.Range("E37").Value = '  +3.42%  '

# Row 38
$ws.Range("B38").Value = 'Stacks'
$ws.Range("C38").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.14'
$ws.Range("E38").Value = '  +26.80%  '

# Row 39
$ws.Range("B39").Value = 'LidoDAOToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.55'
$ws.Range("E39").Value = '  +4.93%  '

# Row 40
$ws.Range("E40").Value = '  -0.01%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '17.63'
$ws.Range("E41").Value = '  +6.74%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '134.44'
$ws.Range("E42").Value = '  +3.00%  '

# Row 43
$ws.Range("E43").Value = '  +5.75%  '

# Row 44
$ws.Range("E44").Value = '  +4.60%  '

# Row 45
$ws.Range("E45").Value = '  +6.11%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.285'
$ws.Range("E46").Value = '  -3.13%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '22.34'
$ws.Range("E47").Value = '  +3.22%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.19'
$ws.Range("E48").Value = '  +56.17%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.149.77'
$ws.Range("E49").Value = '  +4.68%  '

# Row 50
$ws.Range("E50").Value = '  +1.22%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.42'
$ws.Range("E51").Value = '  -4.12%  '
